$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking price strings so they are not
# auto-converted to numbers by Excel type inference (matches source data
# which stores these as plain text/inline strings).
$textCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D12", "D14", "D16", "D20", "D22", "D24", "D27", "D28", "D29", "D32", "D33", "D34", "D35", "D38", "D39", "D41", "D42", "D44", "D45", "D47", "D48", "D49", "D50")
foreach ($tc in $textCells) {
    $ws.Range($tc).NumberFormat = "@"
}

$ws.Range("D2").Value = "43.234.54"
$ws.Range("E2").Value = "  -1.18%  "
$ws.Range("D3").Value = "2.268.93"
$ws.Range("E3").Value = "  -1.90%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "113.25"
$ws.Range("E5").Value = "  +4.69%  "
$ws.Range("D6").Value = "265.11"
$ws.Range("E6").Value = "  -2.25%  "
$ws.Range("D7").Value = "0.619"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("D9").Value = "0.597"
$ws.Range("E9").Value = "  -3.42%  "
$ws.Range("D10").Value = "48.14"
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("D11").Value = "0.0925"
$ws.Range("E11").Value = "  -1.64%  "
$ws.Range("D12").Value = "8.74"
$ws.Range("E12").Value = "  +4.28%  "
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("D14").Value = "15.38"
$ws.Range("E14").Value = "  -2.51%  "
$ws.Range("D15").Value = "2.608.96"
$ws.Range("E15").Value = "  -1.80%  "
$ws.Range("D16").Value = "0.854"
$ws.Range("E16").Value = "  -1.23%  "
$ws.Range("D17").Value = "2.266.17"
$ws.Range("E17").Value = "  -2.07%  "
$ws.Range("D18").Value = "43.132.67"
$ws.Range("E18").Value = "  -1.50%  "
$ws.Range("E19").Value = "  -3.60%  "
$ws.Range("D20").Value = "6.86"
$ws.Range("E20").Value = "  +8.32%  "
$ws.Range("E21").Value = "  -1.80%  "
$ws.Range("D22").Value = "2.42"
$ws.Range("E22").Value = "  -3.36%  "
$ws.Range("E23").Value = "  +3.67%  "
$ws.Range("D24").Value = "230.71"
$ws.Range("E24").Value = "  -1.65%  "
$ws.Range("E25").Value = "  -1.53%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").Value = "11.32"
$ws.Range("E27").Value = "  -0.66%  "
$ws.Range("D28").Value = "3.90"
$ws.Range("E28").Value = "  -1.05%  "
$ws.Range("D29").Value = "40.64"
$ws.Range("E29").Value = "  -3.92%  "
$ws.Range("E30").Value = "  -2.44%  "
$ws.Range("E31").Value = "  -1.77%  "
$ws.Range("D32").Value = "171.07"
$ws.Range("E32").Value = "  -3.86%  "
$ws.Range("D33").Value = "21.23"
$ws.Range("E33").Value = "  -3.41%  "
$ws.Range("D34").Value = "0.0906"
$ws.Range("E34").Value = "  -1.32%  "
$ws.Range("D35").Value = "5.64"
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("E36").Value = "  -0.54%  "
$ws.Range("E37").Value = "  -3.72%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "3.84"
$ws.Range("E38").Value = "  -2.21%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.0350"
$ws.Range("E39").Value = "  -2.44%  "
$ws.Range("E40").Value = "  -7.71%  "
$ws.Range("D41").Value = "14.25"
$ws.Range("E41").Value = "  +15.19%  "
$ws.Range("D42").Value = "75.48"
$ws.Range("E42").Value = "  +12.45%  "
$ws.Range("E43").Value = "  +2.63%  "
$ws.Range("D44").Value = "0.234"
$ws.Range("E44").Value = "  -1.75%  "
$ws.Range("D45").Value = "6.09"
$ws.Range("E45").Value = "  +7.94%  "
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("D47").Value = "1.36"
$ws.Range("E47").Value = "  -2.98%  "
$ws.Range("D48").Value = "8.61"
$ws.Range("E48").Value = "  -2.41%  "
$ws.Range("D49").Value = "0.0983"
$ws.Range("E49").Value = "  -3.88%  "
$ws.Range("D50").Value = "100.46"
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("E51").Value = "  +0.34%  "
